# Applies the "Handles float input without breaking stuff" edit:
#  - Fills in the summary block (rows 10-12) with real Right/Wrong/NotAttempt/
#    Max totals and the computed Total/Max marks string.
#  - Marks each answered question's "Student Ans" cell (column A of each
#    3-block layout) with the student's chosen option and colors it with the
#    correctStyle/incorrectStyle named cell style.
#  - Removes the now-unused 3rd "Student Ans/Correct Ans" block (columns G/H)
#    and the 2nd block (columns D/E) for every question row except the first
#    few that still carry data, shrinking the sheet's used range down to
#    A5:E40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "26/112"

# ---- Drop the 3rd "Student Ans / Correct Ans" block (columns G/H) ------
$ws.Range("G15:H40").Clear()

# ---- Per-question "Student Ans" marks (column A) ------------------------
# row -> (answer text, style)
$answers = @{
    16 = @("Option B", "incorrectStyle")
    17 = @("Option D", "correctStyle")
    18 = @("Option B", "correctStyle")
    21 = @("Option C", "correctStyle")
    25 = @("Option A", "correctStyle")
    27 = @("Option A", "correctStyle")
    33 = @("Option D", "correctStyle")
    34 = @("Option A", "incorrectStyle")
}

foreach ($row in $answers.Keys) {
    $val = $answers[$row][0]
    $style = $answers[$row][1]
    $cell = $ws.Range("A$row")
    $cell.Value = $val
    $cell.Style = $style
}

# Row 16's "2nd block" Student Ans (D16) is also now answered/correct.
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

# ---- Drop the 2nd block (columns D/E) everywhere except rows 16-18 -----
$ws.Range("D19:E40").Clear()
